# Scheduled runner update: refresh market-board derived price/profit figures
# (currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ]) across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR crafting leve sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 180.42857
$ws.Range("I12").Value = 182.6
$ws.Range("K12").Value = 182.6
$ws.Range("M12").Value = -12.59999999999999

$ws.Range("H86").Value = 2991.4443
$ws.Range("I86").Value = 2929.2
$ws.Range("K86").Value = 2929.2
$ws.Range("M86").Value = -1806.2

$ws.Range("H89").Value = 2991.4443
$ws.Range("I89").Value = 2929.2
$ws.Range("K89").Value = 14646
$ws.Range("M89").Value = -9030

$ws.Range("H100").Value = 2995
$ws.Range("I100").Value = 2515.111
$ws.Range("K100").Value = 2515.111
$ws.Range("M100").Value = -1974.111

$ws.Range("H112").Value = 1522
$ws.Range("J112").Value = 1544.12
$ws.Range("L112").Value = 4632.36
$ws.Range("N112").Value = -6848.36

$ws.Range("H132").Value = 4343.0312
$ws.Range("I132").Value = 1810.8948
$ws.Range("K132").Value = 5432.6844
$ws.Range("M132").Value = -2902.6844

$ws.Range("H138").Value = 2354.9167
$ws.Range("J138").Value = 3862.3333
$ws.Range("L138").Value = 11586.9999
$ws.Range("N138").Value = -21866.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3279.5334
$ws.Range("I2").Value = 2991.9167
$ws.Range("K2").Value = 2991.9167
$ws.Range("M2").Value = -2878.9167

$ws.Range("H61").Value = 3707.2307
$ws.Range("I61").Value = 3132.6667
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 3132.6667
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -2920.6667
$ws.Range("N61").Value = -5424

$ws.Range("H102").Value = 1257.5526
$ws.Range("I102").Value = 1228.2858
$ws.Range("K102").Value = 1228.2858
$ws.Range("M102").Value = 393.7141999999999

$ws.Range("H116").Value = 3279.5334
$ws.Range("I116").Value = 2991.9167
$ws.Range("K116").Value = 2991.9167
$ws.Range("M116").Value = -697.9167000000002

$ws.Range("H136").Value = 3707.2307
$ws.Range("I136").Value = 3132.6667
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 9398.000100000001
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -6848.000100000001
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3279.5334
$ws.Range("I3").Value = 2991.9167
$ws.Range("K3").Value = 2991.9167
$ws.Range("M3").Value = -2877.9167

$ws.Range("H82").Value = 19136.637
$ws.Range("I82").Value = 8878.25
$ws.Range("J82").Value = 24998.572
$ws.Range("K82").Value = 8878.25
$ws.Range("L82").Value = 24998.572
$ws.Range("M82").Value = -8495.25
$ws.Range("N82").Value = -25764.572

$ws.Range("H85").Value = 19136.637
$ws.Range("I85").Value = 8878.25
$ws.Range("J85").Value = 24998.572
$ws.Range("K85").Value = 8878.25
$ws.Range("L85").Value = 24998.572
$ws.Range("M85").Value = -7552.25
$ws.Range("N85").Value = -27650.572

$ws.Range("H86").Value = 1978.8889
$ws.Range("I86").Value = 1659.8182
$ws.Range("J86").Value = 3382.8
$ws.Range("K86").Value = 1659.8182
$ws.Range("L86").Value = 3382.8
$ws.Range("M86").Value = -536.8181999999999
$ws.Range("N86").Value = -5628.8

$ws.Range("H89").Value = 1978.8889
$ws.Range("I89").Value = 1659.8182
$ws.Range("J89").Value = 3382.8
$ws.Range("K89").Value = 8299.091
$ws.Range("L89").Value = 16914
$ws.Range("M89").Value = -2683.091
$ws.Range("N89").Value = -28146

$ws.Range("H94").Value = 575.1177
$ws.Range("I94").Value = 589.5
$ws.Range("J94").Value = 508
$ws.Range("K94").Value = 589.5
$ws.Range("L94").Value = 508
$ws.Range("M94").Value = -138.5
$ws.Range("N94").Value = -1410

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2171.4285
$ws.Range("I16").Value = 1866.6666
$ws.Range("K16").Value = 1866.6666
$ws.Range("M16").Value = -1579.6666

$ws.Range("H31").Value = 1688.6
$ws.Range("I31").Value = 1688.6
$ws.Range("K31").Value = 1688.6
$ws.Range("M31").Value = -1393.6

$ws.Range("H34").Value = 1688.6
$ws.Range("I34").Value = 1688.6
$ws.Range("K34").Value = 1688.6
$ws.Range("M34").Value = -1486.6

$ws.Range("H58").Value = 1978.3
$ws.Range("I58").Value = 1656.8
$ws.Range("J58").Value = 2299.8
$ws.Range("K58").Value = 1656.8
$ws.Range("L58").Value = 2299.8
$ws.Range("M58").Value = -1453.8
$ws.Range("N58").Value = -2705.8

$ws.Range("H94").Value = 2867.3333
$ws.Range("I94").Value = 2998
$ws.Range("J94").Value = 2214
$ws.Range("K94").Value = 2998
$ws.Range("L94").Value = 2214
$ws.Range("M94").Value = -2547
$ws.Range("N94").Value = -3116

$ws.Range("H99").Value = 15348
$ws.Range("J99").Value = 20724.8
$ws.Range("L99").Value = 20724.8
$ws.Range("N99").Value = -23720.8

$ws.Range("H113").Value = 2171.4285
$ws.Range("I113").Value = 1866.6666
$ws.Range("K113").Value = 1866.6666
$ws.Range("M113").Value = 303.3334

$ws.Range("H124").Value = 61330.668
$ws.Range("J124").Value = 61330.668
$ws.Range("L124").Value = 61330.668
$ws.Range("N124").Value = -66240.66800000001

$ws.Range("H126").Value = 15348
$ws.Range("J126").Value = 20724.8
$ws.Range("L126").Value = 62174.39999999999
$ws.Range("N126").Value = -67114.39999999999

$ws.Range("H132").Value = 2105.4333
$ws.Range("I132").Value = 2123.9614
$ws.Range("K132").Value = 6371.8842
$ws.Range("M132").Value = -3841.8842

$ws.Range("H134").Value = 2185.2964
$ws.Range("I134").Value = 2080.1428
$ws.Range("K134").Value = 6240.428400000001
$ws.Range("M134").Value = -3705.428400000001

$ws.Range("H136").Value = 1978.3
$ws.Range("I136").Value = 1656.8
$ws.Range("J136").Value = 2299.8
$ws.Range("K136").Value = 4970.4
$ws.Range("L136").Value = 6899.400000000001
$ws.Range("M136").Value = -2420.4
$ws.Range("N136").Value = -11999.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 2700
$ws.Range("I109").Value = 2700
$ws.Range("K109").Value = 8100
$ws.Range("M109").Value = -7060

$ws.Range("H137").Value = 4549128.5
$ws.Range("I137").Value = 8335091
$ws.Range("J137").Value = 5973.8
$ws.Range("K137").Value = 25005273
$ws.Range("L137").Value = 17921.4
$ws.Range("M137").Value = -25000173
$ws.Range("N137").Value = -28121.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9876.4375
$ws.Range("J80").Value = 3986.6667
$ws.Range("L80").Value = 3986.6667
$ws.Range("N80").Value = -5982.6667

$ws.Range("H83").Value = 9876.4375
$ws.Range("J83").Value = 3986.6667
$ws.Range("L83").Value = 19933.3335
$ws.Range("N83").Value = -29917.3335

$ws.Range("H97").Value = 41843.625
$ws.Range("I97").Value = 55058.25
$ws.Range("K97").Value = 55058.25
$ws.Range("M97").Value = -54562.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3059.9333
$ws.Range("J46").Value = 3400
$ws.Range("L46").Value = 3400
$ws.Range("N46").Value = -3776

$ws.Range("H61").Value = 10676.071
$ws.Range("I61").Value = 12541.818
$ws.Range("K61").Value = 12541.818
$ws.Range("M61").Value = -12339.818

$ws.Range("H93").Value = 1631
$ws.Range("I93").Value = 1513.1538
$ws.Range("K93").Value = 1513.1538
$ws.Range("M93").Value = -265.1538

$ws.Range("H113").Value = 10676.071
$ws.Range("I113").Value = 12541.818
$ws.Range("K113").Value = 12541.818
$ws.Range("M113").Value = -10371.818

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3317.2727
$ws.Range("I81").Value = 3784.2856
$ws.Range("K81").Value = 7568.5712
$ws.Range("M81").Value = -6507.5712

$ws.Range("H84").Value = 3317.2727
$ws.Range("I84").Value = 3784.2856
$ws.Range("K84").Value = 37842.856
$ws.Range("M84").Value = -32538.856

$ws.Range("H96").Value = 2123.25
$ws.Range("I96").Value = 1831.5
$ws.Range("J96").Value = 2998.5
$ws.Range("K96").Value = 1831.5
$ws.Range("L96").Value = 2998.5
$ws.Range("M96").Value = -458.5
$ws.Range("N96").Value = -5744.5

$ws.Range("H100").Value = 5362.222
$ws.Range("I100").Value = 8190.25
$ws.Range("K100").Value = 16380.5
$ws.Range("M100").Value = -15839.5

$ws.Range("H132").Value = 8459.096
$ws.Range("I132").Value = 10028.8
$ws.Range("K132").Value = 30086.4
$ws.Range("M132").Value = -27556.4
